# Auto-generated edit script: updates crypto price/volume table cells
# per commit "Updated cryptos list on Mon Sep 18 18:11:47 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.922.46"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.648.49"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.253"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.878.63"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "1.659.66"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").Value = "26.930.82"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.86%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").Value = "1.250.33"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0175"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "1.791.17"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.29%  "
